$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J5").Value = 100
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 27
$ws.Range("N5").Value = -330
$ws.Range("I5").Value = 88
$ws.Range("K5").Value = 88
$ws.Range("H5").Value = 89.5
$ws.Range("M18").Value = -216
$ws.Range("I18").Value = 500
$ws.Range("K18").Value = 500
$ws.Range("H18").Value = 500
$ws.Range("H29").Value = 2542.8572
$ws.Range("M32").Value = -10173.75
$ws.Range("I32").Value = 10499.75
$ws.Range("K32").Value = 10499.75
$ws.Range("H32").Value = 13848.75
$ws.Range("M76").Value = -3380.75
$ws.Range("I76").Value = 3695.75
$ws.Range("K76").Value = 3695.75
$ws.Range("H76").Value = 3956.6
$ws.Range("M79").Value = -2603.75
$ws.Range("I79").Value = 3695.75
$ws.Range("K79").Value = 3695.75
$ws.Range("H79").Value = 3956.6
$ws.Range("M111").Value = 2401.6
$ws.Range("I111").Value = 221.8
$ws.Range("K111").Value = 665.4000000000001
$ws.Range("H111").Value = 468.16666
$ws.Range("J114").Value = 69990
$ws.Range("L114").Value = 69990
$ws.Range("N114").Value = -78668
$ws.Range("H114").Value = 69990
$ws.Range("J121").Value = 1035.1666
$ws.Range("L121").Value = 3105.4998
$ws.Range("N121").Value = -6599.4998
$ws.Range("H121").Value = 1035.1666
$ws.Range("J125").Value = 17859004
$ws.Range("L125").Value = 160731036
$ws.Range("M125").Value = -10171607.1
$ws.Range("N125").Value = -160735956
$ws.Range("I125").Value = 1130451.9
$ws.Range("K125").Value = 10174067.1
$ws.Range("H125").Value = 6453173
$ws.Range("J134").Value = 120000
$ws.Range("L134").Value = 120000
$ws.Range("N134").Value = -130140
$ws.Range("H134").Value = 120000
$ws.Range("M137").Value = -13010.2002
$ws.Range("I137").Value = 5186.7334
$ws.Range("K137").Value = 15560.2002
$ws.Range("H137").Value = 5410
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H140").Value = 0

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M2").Value = -978685.0600000001
$ws.Range("I2").Value = 978798.0600000001
$ws.Range("K2").Value = 978798.0600000001
$ws.Range("H2").Value = 825156.2
$ws.Range("J5").Value = 50.8
$ws.Range("L5").Value = 50.8
$ws.Range("M5").Value = 7.44444
$ws.Range("N5").Value = -274.8
$ws.Range("I5").Value = 104.55556
$ws.Range("K5").Value = 104.55556
$ws.Range("H5").Value = 85.35714
$ws.Range("M31").Value = -4676.5
$ws.Range("I31").Value = 4970.5
$ws.Range("K31").Value = 4970.5
$ws.Range("H31").Value = 4970.5
$ws.Range("M32").Value = -9825.666999999999
$ws.Range("I32").Value = 10112.667
$ws.Range("K32").Value = 10112.667
$ws.Range("H32").Value = 19290.412
$ws.Range("M35").Value = -460.6667
$ws.Range("I35").Value = 866.6667
$ws.Range("K35").Value = 866.6667
$ws.Range("H35").Value = 866.6667
$ws.Range("M61").Value = -60610444
$ws.Range("I61").Value = 60610656
$ws.Range("K61").Value = 60610656
$ws.Range("H61").Value = 58828116
$ws.Range("M74").Value = -19235562
$ws.Range("I74").Value = 19236436
$ws.Range("K74").Value = 19236436
$ws.Range("H74").Value = 18524264
$ws.Range("M77").Value = -96177812
$ws.Range("I77").Value = 19236436
$ws.Range("K77").Value = 96182180
$ws.Range("H77").Value = 18524264
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J88").Value = 1743.5714
$ws.Range("L88").Value = 1743.5714
$ws.Range("N88").Value = -2555.5714
$ws.Range("H88").Value = 1800.625
$ws.Range("J91").Value = 1743.5714
$ws.Range("L91").Value = 1743.5714
$ws.Range("N91").Value = -4551.5714
$ws.Range("H91").Value = 1800.625
$ws.Range("J97").Value = 1499.5
$ws.Range("L97").Value = 1499.5
$ws.Range("M97").Value = -574.2307000000001
$ws.Range("N97").Value = -2491.5
$ws.Range("I97").Value = 1070.2307
$ws.Range("K97").Value = 1070.2307
$ws.Range("H97").Value = 1127.4667
$ws.Range("J110").Value = 2199.6
$ws.Range("L110").Value = 2199.6
$ws.Range("M110").Value = -166296.5
$ws.Range("N110").Value = -6289.6
$ws.Range("I110").Value = 168341.5
$ws.Range("K110").Value = 168341.5
$ws.Range("H110").Value = 92822.45
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H111").Value = 0
$ws.Range("M116").Value = -976504.0600000001
$ws.Range("I116").Value = 978798.0600000001
$ws.Range("K116").Value = 978798.0600000001
$ws.Range("H116").Value = 825156.2
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H117").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -12072.1432
$ws.Range("N122").ClearContents()
$ws.Range("I122").Value = 4840.7144
$ws.Range("K122").Value = 14522.1432
$ws.Range("H122").Value = 4840.7144
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9709016.600000001
$ws.Range("N132").ClearContents()
$ws.Range("I132").Value = 3237182.2
$ws.Range("K132").Value = 9711546.600000001
$ws.Range("H132").Value = 3237182.2
$ws.Range("M136").Value = -181829418
$ws.Range("I136").Value = 60610656
$ws.Range("K136").Value = 181831968
$ws.Range("H136").Value = 58828116

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M3").Value = -978684.0600000001
$ws.Range("I3").Value = 978798.0600000001
$ws.Range("K3").Value = 978798.0600000001
$ws.Range("H3").Value = 825156.2
$ws.Range("J4").Value = 50.8
$ws.Range("L4").Value = 50.8
$ws.Range("M4").Value = 10.44444
$ws.Range("N4").Value = -280.8
$ws.Range("I4").Value = 104.55556
$ws.Range("K4").Value = 104.55556
$ws.Range("H4").Value = 85.35714
$ws.Range("J22").Value = 47619050
$ws.Range("L22").Value = 47619050
$ws.Range("M22").Value = -2235
$ws.Range("N22").Value = -47619396
$ws.Range("I22").Value = 2408
$ws.Range("K22").Value = 2408
$ws.Range("H22").Value = 9525736
$ws.Range("J25").Value = 24983
$ws.Range("L25").Value = 24983
$ws.Range("M25").Value = -389
$ws.Range("N25").Value = -25453
$ws.Range("I25").Value = 624
$ws.Range("K25").Value = 624
$ws.Range("H25").Value = 15239.4
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("J86").Value = 3399.5
$ws.Range("L86").Value = 3399.5
$ws.Range("M86").Value = -2384.4
$ws.Range("N86").Value = -5645.5
$ws.Range("I86").Value = 3507.4
$ws.Range("K86").Value = 3507.4
$ws.Range("H86").Value = 3489.4167
$ws.Range("J89").Value = 3399.5
$ws.Range("L89").Value = 16997.5
$ws.Range("M89").Value = -11921
$ws.Range("N89").Value = -28229.5
$ws.Range("I89").Value = 3507.4
$ws.Range("K89").Value = 17537
$ws.Range("H89").Value = 3489.4167
$ws.Range("M99").Value = -137.5
$ws.Range("I99").Value = 1635.5
$ws.Range("K99").Value = 1635.5
$ws.Range("H99").Value = 1882.5
$ws.Range("M102").ClearContents()
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("H102").Value = 0
$ws.Range("J105").Value = 697
$ws.Range("L105").Value = 697
$ws.Range("M105").Value = 947.4
$ws.Range("N105").Value = -4191
$ws.Range("I105").Value = 799.6
$ws.Range("K105").Value = 799.6
$ws.Range("H105").Value = 782.5
$ws.Range("M134").Value = -46882128
$ws.Range("I134").Value = 15628221
$ws.Range("K134").Value = 46884663
$ws.Range("H134").Value = 15154881

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J22").Value = 832.8
$ws.Range("L22").Value = 832.8
$ws.Range("M22").Value = -16664
$ws.Range("N22").Value = -1532.8
$ws.Range("I22").Value = 17014
$ws.Range("K22").Value = 17014
$ws.Range("H22").Value = 9658.909
$ws.Range("J31").Value = 12999.25
$ws.Range("L31").Value = 12999.25
$ws.Range("M31").Value = -10606.444
$ws.Range("N31").Value = -13589.25
$ws.Range("I31").Value = 10901.444
$ws.Range("K31").Value = 10901.444
$ws.Range("H31").Value = 12244.04
$ws.Range("J34").Value = 12999.25
$ws.Range("L34").Value = 12999.25
$ws.Range("M34").Value = -10699.444
$ws.Range("N34").Value = -13403.25
$ws.Range("I34").Value = 10901.444
$ws.Range("K34").Value = 10901.444
$ws.Range("H34").Value = 12244.04
$ws.Range("M52").Value = -45705.332
$ws.Range("I52").Value = 45999.332
$ws.Range("K52").Value = 45999.332
$ws.Range("H52").Value = 45999.332
$ws.Range("M58").Value = -33342153
$ws.Range("I58").Value = 33342356
$ws.Range("K58").Value = 33342356
$ws.Range("H58").Value = 22734216
$ws.Range("M70").ClearContents()
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J99").Value = 3716.25
$ws.Range("L99").Value = 3716.25
$ws.Range("M99").Value = -2659
$ws.Range("N99").Value = -6712.25
$ws.Range("I99").Value = 4157
$ws.Range("K99").Value = 4157
$ws.Range("H99").Value = 3885.7693
$ws.Range("J105").Value = 3008.6667
$ws.Range("L105").Value = 3008.6667
$ws.Range("M105").Value = -2976222.5
$ws.Range("N105").Value = -6502.6667
$ws.Range("I105").Value = 2977969.5
$ws.Range("K105").Value = 2977969.5
$ws.Range("H105").Value = 1986315.9
$ws.Range("J122").Value = 4045.4546
$ws.Range("L122").Value = 12136.3638
$ws.Range("M122").Value = -10733.353
$ws.Range("N122").Value = -17036.3638
$ws.Range("I122").Value = 4394.451
$ws.Range("K122").Value = 13183.353
$ws.Range("H122").Value = 4332.532
$ws.Range("J126").Value = 3716.25
$ws.Range("L126").Value = 11148.75
$ws.Range("M126").Value = -10001
$ws.Range("N126").Value = -16088.75
$ws.Range("I126").Value = 4157
$ws.Range("K126").Value = 12471
$ws.Range("H126").Value = 3885.7693
$ws.Range("M136").Value = -100024518
$ws.Range("I136").Value = 33342356
$ws.Range("K136").Value = 100027068
$ws.Range("H136").Value = 22734216

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J5").Value = 5752.5
$ws.Range("L5").Value = 17257.5
$ws.Range("N5").Value = -17481.5
$ws.Range("H5").Value = 49036.57
$ws.Range("M44").Value = -346
$ws.Range("I44").Value = 248
$ws.Range("K44").Value = 744
$ws.Range("H44").Value = 286
$ws.Range("M56").Value = -16879.656
$ws.Range("I56").Value = 17409.656
$ws.Range("K56").Value = 17409.656
$ws.Range("H56").Value = 17409.656
$ws.Range("M58").Value = -172
$ws.Range("I58").Value = 100
$ws.Range("K58").Value = 300
$ws.Range("H58").Value = 4274.75
$ws.Range("J68").Value = 713
$ws.Range("L68").Value = 2139
$ws.Range("M68").Value = -770
$ws.Range("N68").Value = -3761
$ws.Range("I68").Value = 527
$ws.Range("K68").Value = 1581
$ws.Range("H68").Value = 564.2
$ws.Range("J69").Value = 12996
$ws.Range("L69").Value = 38988
$ws.Range("M69").Value = -1736
$ws.Range("N69").Value = -40610
$ws.Range("I69").Value = 849
$ws.Range("K69").Value = 2547
$ws.Range("H69").Value = 8947
$ws.Range("J71").Value = 713
$ws.Range("L71").Value = 6417
$ws.Range("M71").Value = -687
$ws.Range("N71").Value = -14529
$ws.Range("I71").Value = 527
$ws.Range("K71").Value = 4743
$ws.Range("H71").Value = 564.2
$ws.Range("J72").Value = 12996
$ws.Range("L72").Value = 116964
$ws.Range("M72").Value = -3585
$ws.Range("N72").Value = -125076
$ws.Range("I72").Value = 849
$ws.Range("K72").Value = 7641
$ws.Range("H72").Value = 8947
$ws.Range("J80").Value = 2999.5
$ws.Range("L80").Value = 8998.5
$ws.Range("M80").Value = -6187.5
$ws.Range("N80").Value = -10870.5
$ws.Range("I80").Value = 2374.5
$ws.Range("K80").Value = 7123.5
$ws.Range("H80").Value = 2687
$ws.Range("J83").Value = 2999.5
$ws.Range("L83").Value = 26995.5
$ws.Range("M83").Value = -16690.5
$ws.Range("N83").Value = -36355.5
$ws.Range("I83").Value = 2374.5
$ws.Range("K83").Value = 21370.5
$ws.Range("H83").Value = 2687
$ws.Range("J92").Value = 999
$ws.Range("L92").Value = 2997
$ws.Range("M92").Value = -167.33331
$ws.Range("N92").Value = -5493
$ws.Range("I92").Value = 471.77777
$ws.Range("K92").Value = 1415.33331
$ws.Range("H92").Value = 524.5
$ws.Range("J135").Value = 5752.5
$ws.Range("L135").Value = 51772.5
$ws.Range("N135").Value = -56842.5
$ws.Range("H135").Value = 49036.57
$ws.Range("J137").Value = 3887.3333
$ws.Range("L137").Value = 11661.9999
$ws.Range("N137").Value = -21861.9999
$ws.Range("H137").Value = 7146873

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J15").Value = 29589.166
$ws.Range("L15").Value = 29589.166
$ws.Range("N15").Value = -30165.166
$ws.Range("H15").Value = 28362.143
$ws.Range("J20").Value = 33953
$ws.Range("L20").Value = 33953
$ws.Range("M20").Value = -5002257.5
$ws.Range("N20").Value = -34443
$ws.Range("I20").Value = 5002502.5
$ws.Range("K20").Value = 5002502.5
$ws.Range("H20").Value = 2518227.8
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H24").Value = 7511250
$ws.Range("J81").Value = 29589.166
$ws.Range("L81").Value = 29589.166
$ws.Range("N81").Value = -31585.166
$ws.Range("H81").Value = 28362.143
$ws.Range("J84").Value = 29589.166
$ws.Range("L84").Value = 88767.49800000001
$ws.Range("N84").Value = -98751.49800000001
$ws.Range("H84").Value = 28362.143
$ws.Range("M102").Value = 549.25
$ws.Range("I102").Value = 1072.75
$ws.Range("K102").Value = 1072.75
$ws.Range("H102").Value = 1066.3334
$ws.Range("J103").Value = 45000
$ws.Range("L103").Value = 45000
$ws.Range("N103").Value = -47344
$ws.Range("H103").Value = 45000
$ws.Range("J122").Value = 2300
$ws.Range("L122").Value = 6900
$ws.Range("M122").Value = -457240.36
$ws.Range("N122").Value = -11800
$ws.Range("I122").Value = 153230.12
$ws.Range("K122").Value = 459690.36
$ws.Range("H122").Value = 136460.11
$ws.Range("J132").Value = 5961.3335
$ws.Range("L132").Value = 17884.0005
$ws.Range("N132").Value = -22944.0005
$ws.Range("H132").Value = 1988958.2

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J16").Value = 4708.909
$ws.Range("L16").Value = 4708.909
$ws.Range("M16").Value = -1777.25
$ws.Range("N16").Value = -5048.909
$ws.Range("I16").Value = 1947.25
$ws.Range("K16").Value = 1947.25
$ws.Range("H16").Value = 3972.4666
$ws.Range("J46").Value = 2494
$ws.Range("L46").Value = 2494
$ws.Range("M46").Value = -1211.4
$ws.Range("N46").Value = -2870
$ws.Range("I46").Value = 1399.4
$ws.Range("K46").Value = 1399.4
$ws.Range("H46").Value = 1581.8334
$ws.Range("J61").Value = 4161.5
$ws.Range("L61").Value = 4161.5
$ws.Range("M61").Value = -3811.3547
$ws.Range("N61").Value = -4565.5
$ws.Range("I61").Value = 4013.3547
$ws.Range("K61").Value = 4013.3547
$ws.Range("H61").Value = 4037.3784
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H62").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("M82").Value = -1145
$ws.Range("I82").Value = 1506
$ws.Range("K82").Value = 1506
$ws.Range("H82").Value = 1643.75
$ws.Range("M85").Value = -258
$ws.Range("I85").Value = 1506
$ws.Range("K85").Value = 1506
$ws.Range("H85").Value = 1643.75
$ws.Range("M99").ClearContents()
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("H99").Value = 39999
$ws.Range("M107").Value = -2637
$ws.Range("I107").Value = 4557
$ws.Range("K107").Value = 4557
$ws.Range("H107").Value = 4557
$ws.Range("J113").Value = 4161.5
$ws.Range("L113").Value = 4161.5
$ws.Range("M113").Value = -1843.3547
$ws.Range("N113").Value = -8501.5
$ws.Range("I113").Value = 4013.3547
$ws.Range("K113").Value = 4013.3547
$ws.Range("H113").Value = 4037.3784
$ws.Range("J122").Value = 6750
$ws.Range("L122").Value = 20250
$ws.Range("M122").Value = -17254
$ws.Range("N122").Value = -25150
$ws.Range("I122").Value = 6568
$ws.Range("K122").Value = 19704
$ws.Range("H122").Value = 6637.3335
$ws.Range("M136").Value = -5943.75
$ws.Range("I136").Value = 2831.25
$ws.Range("K136").Value = 8493.75
$ws.Range("H136").Value = 2873.1875

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J113").Value = 676.5
$ws.Range("L113").Value = 2029.5
$ws.Range("M113").Value = 1419.00001
$ws.Range("N113").Value = -6369.5
$ws.Range("I113").Value = 250.33333
$ws.Range("K113").Value = 750.99999
$ws.Range("H113").Value = 287.3913
$ws.Range("J132").Value = 13087.777
$ws.Range("L132").Value = 39263.331
$ws.Range("M132").Value = -78957926
$ws.Range("N132").Value = -44323.331
$ws.Range("I132").Value = 26320152
$ws.Range("K132").Value = 78960456
$ws.Range("H132").Value = 17864310
$ws.Range("J136").Value = 1620.4
$ws.Range("L136").Value = 4861.200000000001
$ws.Range("N136").Value = -9961.200000000001
$ws.Range("H136").Value = 21741226
$ws.Range("J138").Value = 139158.5
$ws.Range("L138").Value = 139158.5
$ws.Range("N138").Value = -149438.5
$ws.Range("H138").Value = 139158.5
